# Adapt column header formatting to respective input file names.
#
# The sheet used to contrast two AHB format versions with generic
# "_old" / "_new" suffixed headers (plus a "diff" column in between).
# This renames those headers to carry the actual format-version labels
# ("_FV2310" for the old/left side, "_FV2404" for the new/right side),
# turns the header+data range into a real Excel Table ("Table1") with
# autofilter, and freezes the header row so it stays visible on scroll.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename header labels -----------------------------------
$oldSuffix = "_old"
$newSuffix = "_new"
$fvOldSuffix = "_FV2310"
$fvNewSuffix = "_FV2404"

$headerColumnCount = 21
for ($c = 1; $c -le $headerColumnCount; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        if ($val.EndsWith($oldSuffix)) {
            $base = $val.Substring(0, $val.Length - $oldSuffix.Length)
            $cell.Value = $base + $fvOldSuffix
        } elseif ($val.EndsWith($newSuffix)) {
            $base = $val.Substring(0, $val.Length - $newSuffix.Length)
            $cell.Value = $base + $fvNewSuffix
        }
    }
}

# --- Step 2: convert the range into an Excel Table -------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$tableRange = $ws.Range("A1:U" + $lastRow)

$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- Step 3: freeze the header row -----------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
